$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. Sheet1 (root_hospital_cities): add lat/lng columns with values
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("root_hospital_cities")

$ws1.Range("B1").Value = "lat"
$ws1.Range("C1").Value = "lng"
$ws1.Range("B1").Font.Color = 0
$ws1.Range("C1").Font.Color = 0

$ws1.Range("B2").Font.Color = 0
$ws1.Range("C2").Font.Color = 0
$ws1.Range("B3").Font.Color = 0
$ws1.Range("C3").Font.Color = 0

$ws1.Range("B2").NumberFormat = "0.000000"
$ws1.Range("C2").NumberFormat = "0.000000"
$ws1.Range("B3").NumberFormat = "0.000000"
$ws1.Range("C3").NumberFormat = "0.000000"

$ws1.Range("B2").Value = 40.712784
$ws1.Range("C2").Value = -74.005941
$ws1.Range("B3").Value = 37.151165
$ws1.Range("C3").Value = -88.731998

# ----------------------------------------------------------------------
# 2. Sheet4 (attributes): insert 2 new rows describing lat/lng attrs
# ----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("attributes")
$null = $ws4.Activate()

$null = $ws4.Rows.Item(3).Insert()
$null = $ws4.Rows.Item(3).Insert()

$ws4.Range("A3").Value = "lat"
$ws4.Range("B3").Value = "root_hospital_cities"
$ws4.Range("G3").Value = "latitude in degrees"

$ws4.Range("A4").Value = "lng"
$ws4.Range("B4").Value = "root_hospital_cities"
$ws4.Range("G4").Value = "longitude in degrees"

$ws4.Range("C3").Value = "decimal"
$ws4.Range("C4").Value = "decimal"

$null = $ws4.Range("C5").Select()

# ----------------------------------------------------------------------
# 3. Re-activate sheet1 (root_hospital_cities) as the active tab, with
#    selection returned to the default top-left cell.
# ----------------------------------------------------------------------
$null = $ws1.Activate()
$null = $ws1.Range("A1").Select()

Write-Output "done"
